$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:numero-de-personas-en-el-hogar"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "null"

# Row 3: dim/medida
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "null"

# Row 4: type/URI
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Provincia"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-comarca"
$ws.Range("H4").Value = "null"

# Row 5 no longer exists - remove it entirely
$ws.Rows.Item(5).Delete()
